$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Merge the "WED Sep 27" + " 18:39:51 PDT 2017" runs into one run
#    (Find/Replace naturally collapses the two runs into a single
#    run when the replacement text spans both of them).
# -----------------------------------------------------------------
$null = $d.Content.Find.Execute("WED Sep 27 18:39:51 PDT 2017", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "WED Sep 27 18:39:51 PDT 2017", 2)

# -----------------------------------------------------------------
# 2) Bold the last "Amount balance ... - 274957.0" paragraph
#    (paragraph mark + every run in it).
# -----------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("- 274957.0")
$amtPara = $rng.Paragraphs(1)
$amtPara.Range.Font.Bold = 1

# -----------------------------------------------------------------
# 3) Insert the new 29/09/2017 purchase record right after that
#    paragraph, followed by a bold blank line and a plain blank
#    line (matching the pattern used after every other record).
# -----------------------------------------------------------------
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$rPrPlain = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr>'
$rPrBold  = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/></w:rPr>'
$pPrPlain = '<w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr></w:pPr>'
$pPrBold  = '<w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/></w:rPr></w:pPr>'

$xml = "<w:p $w>$pPrBold</w:p>" +
       "<w:p $w>$pPrPlain<w:r>$rPrPlain<w:t>THU Sep 28</w:t></w:r><w:r>$rPrPlain<w:t xml:space=`"preserve`"> 14:05:43 PDT 2017</w:t></w:r></w:p>" +
       "<w:p $w>$pPrPlain<w:r>$rPrPlain<w:t>Person Name</w:t></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/><w:t>- NG</w:t></w:r></w:p>" +
       "<w:p $w>$pPrPlain<w:r>$rPrPlain<w:t>---------------------------------------------------------------</w:t></w:r></w:p>" +
       "<w:p $w>$pPrPlain<w:r>$rPrPlain<w:t>Item Name</w:t></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/><w:t>- CARROT</w:t></w:r></w:p>" +
       "<w:p $w>$pPrPlain<w:r>$rPrPlain<w:t>Number of Pockets</w:t></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/><w:t>- 5</w:t></w:r></w:p>" +
       "<w:p $w>$pPrPlain<w:r>$rPrPlain<w:t>Number of KGs</w:t></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/><w:t>- 464</w:t></w:r></w:p>" +
       "<w:p $w>$pPrPlain<w:r>$rPrPlain<w:t>Rate</w:t></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/><w:t>- 20</w:t></w:r></w:p>" +
       "<w:p $w>$pPrPlain<w:r>$rPrPlain<w:t>Transport &amp; Miscellaneous</w:t></w:r><w:r>$rPrPlain<w:tab/><w:t>- 75</w:t></w:r></w:p>" +
       "<w:p $w>$pPrPlain<w:r>$rPrPlain<w:t>Total Price</w:t></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/></w:r><w:r>$rPrPlain<w:tab/><w:t>- 9355.0</w:t></w:r></w:p>" +
       "<w:p $w>$pPrBold<w:r>$rPrBold<w:t>Amount balance</w:t></w:r><w:r>$rPrBold<w:tab/></w:r><w:r>$rPrBold<w:tab/></w:r><w:r>$rPrBold<w:tab/><w:t>- 284312.0</w:t></w:r></w:p>" +
       "<w:p $w>$pPrBold</w:p>" +
       "<w:p $w>$pPrPlain</w:p>"

$nextPara = $amtPara.Next()
$insertPoint = $nextPara.Range
$insertPoint.Collapse(1)
$insertPoint.InsertXML($xml)

Write-Output "done"
